$d = $word.ActiveDocument

# 1. "英文" -> "英语" (2 occurrences: hyperlink run + standalone paragraph run)
$d.Content.Find.Execute("英文", $true, $false, $false, $false, $false, $true, 1, $false, "英语", 2) | Out-Null

# 2. Language list red text
$d.Content.Find.Execute(" / 葡萄牙文 / 法文 / 泰文 / 越南文 / 西班牙文", $true, $false, $false, $false, $false, $true, 1, $false, " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语", 2) | Out-Null

# 3. Table cell: "簡介" -> "简介"
$d.Content.Find.Execute("簡介", $true, $false, $false, $false, $false, $true, 1, $false, "简介", 2) | Out-Null

# 4. Table cell description paragraph
$d.Content.Find.Execute("寄給參加活動的合作夥伴的電子郵件。 此電子郵件將包括照片畫廊，將通過 customer.io 發送。", $true, $false, $false, $false, $false, $true, 1, $false, "一封发送给参加活动的合作伙伴的邮件。 这封邮件将包含一个照片画廊，将通过 customer.io 发送。", 2) | Out-Null

# 5. "目標受眾" -> "目标受众"
$d.Content.Find.Execute("目標受眾", $true, $false, $false, $false, $false, $true, 1, $false, "目标受众", 2) | Out-Null

# 6. "活動參加者" -> "活动参与者"
$d.Content.Find.Execute("活動參加者", $true, $false, $false, $false, $false, $true, 1, $false, "活动参与者", 2) | Out-Null

# 7. "Subject: " -> "主题: "
$d.Content.Find.Execute("Subject: ", $true, $false, $false, $false, $false, $true, 1, $false, "主题: ", 2) | Out-Null

# 8. "感謝您參加 " -> "感谢您参加 " (2 occurrences: subject line + greeting paragraph)
$d.Content.Find.Execute("感謝您參加 ", $true, $false, $false, $false, $false, $true, 1, $false, "感谢您参加 ", 2) | Out-Null

# 9. Heading paragraph
$d.Content.Find.Execute("您使我們的活動圓滿成功！ 🎉", $true, $false, $false, $false, $false, $true, 1, $false, "您使我们的活动圆满成功！ 🎉", 2) | Out-Null

# 10. "[PARTNER NAME]" -> "[合作伙伴姓名]"
$d.Content.Find.Execute("[PARTNER NAME]", $true, $false, $false, $false, $false, $true, 1, $false, "[合作伙伴姓名]", 2) | Out-Null

# 11. "， " (after partner name placeholder) -> ", "
$d.Content.Find.Execute("， ", $true, $false, $false, $false, $false, $true, 1, $false, ", ", 2) | Out-Null

# 12. " 於 " -> " 于 "
$d.Content.Find.Execute(" 於 ", $true, $false, $false, $false, $false, $true, 1, $false, " 于 ", 2) | Out-Null

# 13. Lone "，" between [CITY] and [COUNTRY] -> ", " (scoped replace via precise Range,
#     since a bare "，" also occurs elsewhere in the same paragraph with different target text)
$city_para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text
    if ($t -ne $null -and $t.Length -gt 0 -and $t.Substring(0,1) -eq [char]24863 -and $t.Contains("[CITY]")) {
        $city_para = $para
        break
    }
}
if ($city_para -ne $null) {
    $ptext = $city_para.Range.Text
    $pstart = $city_para.Range.Start
    $cityIdx = $ptext.IndexOf("[CITY]")
    $commaIdx = $ptext.IndexOf([char]65292, $cityIdx)
    if ($commaIdx -ge 0) {
        $target = $pstart + $commaIdx
        $rr = $d.Range($target, $target + 1)
        $rr.Text = ", "
    }
}

# 14. Closing sentence in thanks paragraph
$d.Content.Find.Execute("。 希望您度過了愉快的時光，很高興認識您！", $true, $false, $false, $false, $false, $true, 1, $false, "。 希望您玩得开心，很高兴认识您！", 2) | Out-Null

# 15/16. "如需瀏覽會議/" + "研討會/聯盟之旅" -> "如需浏览 " + "会议/研讨会/联盟之旅"
$d.Content.Find.Execute("如需瀏覽會議/", $true, $false, $false, $false, $false, $true, 1, $false, "如需浏览 ", 2) | Out-Null
$d.Content.Find.Execute("研討會/聯盟之旅", $true, $false, $false, $false, $false, $true, 1, $false, "会议/研讨会/联盟之旅", 2) | Out-Null

# 17. Trailing segment after comment reference
$d.Content.Find.Execute("的照片和精彩片段，並隨時了解我們為您舉辦的最新活動和計劃，請關注我們的社交媒體帳戶：", $true, $false, $false, $false, $false, $true, 1, $false, " 的照片和精彩片段，并随时了解我们举办的最新活动和计划，请关注我们：", 2) | Out-Null

# 18. Final paragraph
$d.Content.Find.Execute("希望這次活動能給您們帶來和我們一樣的啟發，讓我們繼續共同成長！", $true, $false, $false, $false, $false, $true, 1, $false, "希望这次活动能给您们带来和我们一样的启发，让我们继续共同成长！", 2) | Out-Null

# 19. Comment text: "選擇其中一個" -> "选择其中之一"
$comment = $d.Comments.Item(1)
$comment.Range.Text = "选择其中之一"
